# Shopzier_CheckOut_Divyansh.xlsx update:
#   - Rename the "Test Cases" sheet to "CheckoutTCs"
#   - Update the selection on the "Scenario" sheet to A3 (was C7)
#   - Keep "CheckoutTCs" as the active/selected sheet/tab

$wb = $excel.ActiveWorkbook

# Rename "Test Cases" -> "CheckoutTCs"
$tcSheet = $wb.Worksheets.Item("Test Cases")
$tcSheet.Name = "CheckoutTCs"

# Update the selected cell on the "Scenario" sheet from C7 to A3.
# Briefly activate Scenario to move its selection, then re-activate
# the CheckoutTCs sheet so it remains the workbook's active tab.
$scenarioSheet = $wb.Worksheets.Item("Scenario")
$scenarioSheet.Activate()
$scenarioSheet.Range("A3").Select()

$tcSheet.Activate()
